# "after adding 2nd review ppt"
#
# 1. Slide 3 ("Base Paper Details"): the author names paragraph gets split
#    into several runs (so that individual names like "Vansh", "Kedia",
#    "Mayand" and "Bickey" can be flagged/marked separately by the
#    spell-checker). The visible text is unchanged.
# 2. Slide 9 ("Modules"): the progress table is resized/stretched to fill
#    the slide, and the "Sprint N" cells are renamed to "Module N".

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 3: split the authors paragraph into multiple runs without
# changing the visible text. Re-asserting the (unchanged) Bold value on
# each sub-range forces PowerPoint to split a fresh run for that span.
# ---------------------------------------------------------------------
$slide3 = $p.Slides.Item(3)
$authorsShape = $slide3.Shapes.Item(2)
$authorsRange = $authorsShape.TextFrame.TextRange

$authorsRange.Characters(1, 13).Font.Bold = $false    # "Aman Bhatia, "
$authorsRange.Characters(14, 5).Font.Bold = $false    # "Vansh"
$authorsRange.Characters(19, 1).Font.Bold = $false    # " "
$authorsRange.Characters(20, 5).Font.Bold = $false    # "Kedia"
$authorsRange.Characters(25, 17).Font.Bold = $false   # ", Anshul Shroff, "
$authorsRange.Characters(42, 6).Font.Bold = $false    # "Mayand"
$authorsRange.Characters(48, 8).Font.Bold = $false    # " Kumar, "
$authorsRange.Characters(56, 6).Font.Bold = $false    # "Bickey"
$authorsRange.Characters(62, 20).Font.Bold = $false   # " Kumar Shah, Aryan ,"

# ---------------------------------------------------------------------
# Slide 9: resize the modules/sprints table and rename "Sprint N" -> "Module N"
# ---------------------------------------------------------------------
$slide9 = $p.Slides.Item(9)
$tableShape = $slide9.Shapes.Item(2)

$tableShape.Width = 806.09455
$tableShape.Height = 339.0705

$tbl = $tableShape.Table
$tbl.Cell(2, 1).Shape.TextFrame.TextRange.Text = "Module 1"
$tbl.Cell(3, 1).Shape.TextFrame.TextRange.Text = "Module 2"
$tbl.Cell(4, 1).Shape.TextFrame.TextRange.Text = "Module 3"
$tbl.Cell(5, 1).Shape.TextFrame.TextRange.Text = "Module 4"
$tbl.Cell(6, 1).Shape.TextFrame.TextRange.Text = "Module 5"
$tbl.Cell(7, 1).Shape.TextFrame.TextRange.Text = "Module 6"
